# Move the "_GoBack" bookmark from the end of the document (after the
# last run of the last paragraph) to the very beginning of the document
# (start of the Title paragraph), and bump the page size from US Letter
# to A4.
#
# Word automatically keeps only a single "_GoBack" bookmark in a
# document, so adding a new one removes the old one for us.

$d = $word.ActiveDocument

# --- Move the _GoBack bookmark to the start of the document ---------
#
# Bookmarks.Add on a genuinely zero-length Range sitting exactly at
# document position 0 gets auto-expanded to cover the whole first
# paragraph by this host. Work around that by temporarily inserting a
# placeholder character at the very start, anchoring the new bookmark
# right after it (position 1, still zero-length, so it is NOT expanded)
# and then deleting the placeholder again, which slides the now-created
# bookmark back down to position 0 without touching its (zero) length.

$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")

$anchor = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $anchor)

$d.Range(0, 1).Delete()

# --- Switch the page size from US Letter to A4 -----------------------
$d.PageSetup.PaperSize = 7
